# Commit: swap the deck's theme colours from the custom "Integral" palette
# over to the stock "Office Theme" palette (Design tab -> Themes -> "Office
# Theme", the first/default built-in theme). Font scheme and format scheme
# (fills/lines/effects) are already identical between the two themes, so
# only the 12 theme colour slots need to change.
#
# PowerPoint's ThemeColorScheme / ColorScheme.Colors collection walks the
# slots in OOXML clrScheme order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2
#   7 accent3 8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
#
# RGB() packs a hex colour RRGGBB the same way VBA's RGB(r,g,b) does:
# value = r + g*256 + b*65536.
function Hex-ToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock "Office Theme" colour scheme.
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $colorScheme.Colors($i).RGB = Hex-ToRgbLong $officeThemeHex[$i - 1]
}
